$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column D ("hours dedicated" notes) for the weekly tracking sheet.
# Order below matches the original authoring order so that the shared
# string table is built up the same way as in the source workbook.
$ws.Range("D7").Value  = '3h (guía de instalación, comunicación con otros grupos)'
$ws.Range("D6").Value  = '8h(reunion, diagramas, readme, modulo html detalles, método get de detalles, funcionalidad de busquedas)'
$ws.Range("D9").Value  = '1h(Documentación)'
$ws.Range("D2").Value  = '1''5h(organización reuniones, documentacion)'
$ws.Range("D10").Value = '2h (reunión, feedback guía de instalación)'
$ws.Range("D4").Value  = '6h(reunión,implementación web y BBDD)'
$ws.Range("D5").Value  = '2h(documentacion)'
$ws.Range("D8").Value  = '1''5h(organización reuniones, documentacion)'

# Match the resulting active-cell selection left behind in the sheet.
$null = $ws.Range("D10").Select()
